# Fruta / hortaliza, semanal
# Updates weekly price/volume data rows (2-21) on the "Alcachofa" sheet:
# dates, variety, volume, min/max/avg prices, unit of sale, origin province,
# price per kg and kg-or-units basis are refreshed to the new weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 12500
$ws.Range("P2").Value = 417

$ws.Range("D3").Value = 44427
$ws.Range("H3").Value = 'Madrigal'
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12500
$ws.Range("N3").Value = '$/caja 40 unidades'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 312
$ws.Range("Q3").Value = 40

$ws.Range("D4").Value = 44498
$ws.Range("H4").Value = 'Española'
$ws.Range("K4").Value = 8500
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8750
$ws.Range("N4").Value = '$/caja 30 unidades'
$ws.Range("O4").Value = 'Provincia de Limarí'
$ws.Range("P4").Value = 292
$ws.Range("Q4").Value = 30

$ws.Range("D5").Value = 44855
$ws.Range("J5").Value = 540
$ws.Range("K5").Value = 9500
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9750
$ws.Range("O5").Value = 'Provincia del Elquí'
$ws.Range("P5").Value = 325

$ws.Range("H6").Value = 'Española'
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 11500
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11750
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("P6").Value = 392
$ws.Range("Q6").Value = 30

$ws.Range("D7").Value = 44426
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 12500
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 12750
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 319

$ws.Range("D8").Value = 45037
$ws.Range("H8").Value = 'Madrigal'
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = '$/caja 40 unidades'
$ws.Range("P8").Value = 412
$ws.Range("Q8").Value = 40

$ws.Range("D9").Value = 44687
$ws.Range("H9").Value = 'Española'
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 19000
$ws.Range("M9").Value = 18500
$ws.Range("N9").Value = '$/caja 30 unidades'
$ws.Range("P9").Value = 617
$ws.Range("Q9").Value = 30

$ws.Range("D10").Value = 44784
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("J10").Value = 520
$ws.Range("K10").Value = 11500
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11750
$ws.Range("N10").Value = '$/caja 40 unidades'
$ws.Range("O10").Value = 'Provincia del Elquí'
$ws.Range("P10").Value = 294
$ws.Range("Q10").Value = 40

$ws.Range("D11").Value = 45049
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 438

$ws.Range("D12").Value = 44420
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 362

$ws.Range("D13").Value = 44420
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("J13").Value = 700
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13500
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 338
$ws.Range("Q13").Value = 40

$ws.Range("D14").Value = 44858
$ws.Range("H14").Value = 'Española'
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 9500
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9750
$ws.Range("N14").Value = '$/caja 30 unidades'
$ws.Range("P14").Value = 325
$ws.Range("Q14").Value = 30

$ws.Range("D15").Value = 44438
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 383

$ws.Range("D16").Value = 44790
$ws.Range("H16").Value = 'Española'
$ws.Range("J16").Value = 400
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("P16").Value = 483
$ws.Range("Q16").Value = 30

$ws.Range("D17").Value = 44790
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 11500
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11750
$ws.Range("P17").Value = 294

$ws.Range("D18").Value = 44701
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("O18").Value = 'Provincia del Elquí'
$ws.Range("P18").Value = 650

$ws.Range("D19").Value = 44767
$ws.Range("H19").Value = 'Madrigal'
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("N19").Value = '$/caja 40 unidades'
$ws.Range("O19").Value = 'Provincia de Limarí'
$ws.Range("P19").Value = 362
$ws.Range("Q19").Value = 40

$ws.Range("D20").Value = 44484
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9500
$ws.Range("O20").Value = 'Provincia del Elquí'
$ws.Range("P20").Value = 317

$ws.Range("D21").Value = 44729
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16500
$ws.Range("P21").Value = 412
